# Restriction list maintenance — add newly-restricted brands
# ("expanded restricted brands in restriction list")
#
# The sheet flags recently-added brand names with a highlighted cell style.
# Previously-flagged rows lose that highlight (they are no longer "new"),
# and the three newly restricted brands (angel, tresor, idole) receive it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

$xlPasteFormats = -4122

# --- capture the two "highlight" looks before they get overwritten -------
# style used by the most recent batch of new entries (full border + fill)
$ws.Range("A6").Copy() | Out-Null
$ws.Range("Z1").PasteSpecial($xlPasteFormats) | Out-Null
# style used by the single trailing/new entry (left+right border + fill)
$ws.Range("A127").Copy() | Out-Null
$ws.Range("Z2").PasteSpecial($xlPasteFormats) | Out-Null

# --- add the newly restricted brands --------------------------------------
$ws.Range("A131").Value = "angel"
$ws.Range("A132").Value = "tresor"
$ws.Range("A133").Value = "idole"

# give the new entries the "new" highlight style (copied from the old batch)
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("A131:A133").PasteSpecial($xlPasteFormats) | Out-Null

# trailing blank row keeps the list's "new-entry" marker going forward
$ws.Range("Z2").Copy() | Out-Null
$ws.Range("A134").PasteSpecial($xlPasteFormats) | Out-Null

# --- the previous "new" rows are no longer new: drop their highlight ------
$ws.Range("A3").Copy() | Out-Null
$plainTargets = "A6","A20","A51","A54","A58","A59","A60","A62","A88","A96","A118","A119","A127"
foreach ($addr in $plainTargets) {
    $ws.Range($addr).PasteSpecial($xlPasteFormats) | Out-Null
}

# --- clean up scratch cells used to stash formats --------------------------
$ws.Range("Z1:Z2").Clear() | Out-Null

$excel.CutCopyMode = $false

# --- update selection / scroll state to land on the newly added rows ------
$ws.Range("A127:A133").Select() | Out-Null
